$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$files = @(
    "SCRIPT/P02P01A/um1203.ssb",
    "SCRIPT/P02P01A/um1206.ssb",
    "SCRIPT/P02P01A/um1302.ssb",
    "SCRIPT/P02P01A/um1305.ssb",
    "SCRIPT/P02P01A/um1308.ssb",
    "SCRIPT/P02P01A/um1403.ssb",
    "SCRIPT/P02P01A/um1406.ssb",
    "SCRIPT/P02P01A/um1503.ssb",
    "SCRIPT/P02P01A/um1603.ssb",
    "SCRIPT/P02P01A/um1606.ssb"
)

$row = 7
foreach ($f in $files) {
    $ws.Cells.Item($row, 1).Value = $f
    $ws.Rows.Item($row).RowHeight = 43.2
    $row = $row + 1
}

$ws.Range("C2").Select()
